# Weekly update: insert a new week's price record for
# "Vega Monumental Concepción - Berenjena" at the top of the data block
# (row 88), pushing the existing records (rows 88-107) down by one row
# to rows 89-108.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 88; existing rows 88..107 shift to 89..108.
$ws.Rows.Item(88).Insert()

# Populate the new row 88 with this week's record.
$ws.Range("A88").Value = 11
$ws.Range("B88").Value = "Vega Monumental Concepción"
$ws.Range("C88").Value = "Bíobío"
$ws.Range("D88").Value = 44782
$ws.Range("E88").Value = 8
$ws.Range("F88").Value = 100112001
$ws.Range("G88").Value = "Berenjena"
$ws.Range("H88").Value = "Sin especificar"
$ws.Range("I88").Value = "Primera"
$ws.Range("J88").Value = 100
$ws.Range("K88").Value = 12000
$ws.Range("L88").Value = 13000
$ws.Range("M88").Value = 12500
$ws.Range("N88").Value = "$/caja 60 unidades"
$ws.Range("O88").Value = "Región de Arica y Parinacota"
$ws.Range("P88").Value = 208
$ws.Range("Q88").Value = 60
$ws.Range("R88").Value = "Hortaliza"
